# edit.ps1 - applies the WBT101 signoff-sheet changes described by the diff:
#  1. Move the "_GoBack" bookmark from inside the 4B.7 "Advanced" cell to a
#     zero-length bookmark immediately before the "Signatures" heading run.
#  2. Merge the split "Adva" / "nced" runs in the 4B.7 Category cell into a
#     single "Advanced" run (this also removes the old bookmark there).
#  3. Split "04D (BLE Centrals)" into three runs: "04D (BLE Central", "s", ")"
#  4. Fill in several previously-empty Category/Description cells under the
#     "04D (BLE Centrals)" section (4D.1 - 4D.5).
#  5. Insert a brand new row (4D.6, Advanced, "Run the Advertising Scanner")
#     after the 4D.5 row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-EmptyCellText {
    # Fills a run-less (truly empty) cell paragraph with a single run of text,
    # using the same color/size/size-complex-script formatting used throughout
    # this table (color 000000, sz 16 half-points = 8pt, szCs 18 = 9pt).
    param($cell, [string]$text)
    $cell.Range.Text = $text
    $rng = $cell.Range
    $rng.Font.Color = 0
    $rng.Font.Size = 8
    $rng.Font.SizeBi = 9
}

function Split-Run {
    # Forces a run boundary at the given document position by toggling Bold
    # on/off on the trailing sub-range; this does not alter visible formatting
    # (it ends up identical to neighboring runs) but prevents the engine from
    # silently re-coalescing adjacent runs that would otherwise look the same.
    param($doc, [int]$pos, [int]$len)
    $r = $doc.Range($pos, $pos + $len)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark to just before "Signatures".
# ---------------------------------------------------------------------------
$startRng = $d.Range(0, 0)
$startRng.InsertBefore("X")
$placeholder = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------------
# 2. Merge "Adva" + "nced" (+ remove the old bookmark) into a single run.
# ---------------------------------------------------------------------------
$advCell = $t.Rows.Item(35).Cells.Item(4)
$advRng = $advCell.Range
$mergedRng = $d.Range($advRng.Start, $advRng.End)
$mergedRng.Text = "Advanced"

# ---------------------------------------------------------------------------
# 3. Split "04D (BLE Centrals)" into "04D (BLE Central" / "s" / ")"
# ---------------------------------------------------------------------------
$bleCell = $t.Rows.Item(43).Cells.Item(2)
$bleStart = $bleCell.Range.Start
Split-Run $d ($bleStart + 16) 1   # "s"
Split-Run $d ($bleStart + 17) 1   # ")"

# ---------------------------------------------------------------------------
# 4. Fill in the previously-empty 4D.1 - 4D.5 cells.
# ---------------------------------------------------------------------------

# 4D.1 (row 43) - Description
Set-EmptyCellText $t.Rows.Item(43).Cells.Item(5) "Make an Observer"

# 4D.2 (row 44) - Description
Set-EmptyCellText $t.Rows.Item(44).Cells.Item(5) "Add a Filter to Show Only Your Device"

# 4D.3 (row 45) - Category + Description (3 runs)
Set-EmptyCellText $t.Rows.Item(45).Cells.Item(4) "Basic"
$descCell = $t.Rows.Item(45).Cells.Item(5)
Set-EmptyCellText $descCell "Connect to Your Peripheral and Turn ON/OFF the LED"
$descStart = $descCell.Range.Start
Split-Run $d ($descStart + 8) 18   # "to Your Peripheral"
Split-Run $d ($descStart + 26) 24  # " and Turn ON/OFF the LED"

# 4D.4 (row 46) - Category + Description
Set-EmptyCellText $t.Rows.Item(46).Cells.Item(4) "Advanced"
Set-EmptyCellText $t.Rows.Item(46).Cells.Item(5) "Add Commands to Turn Notify ON/OFF"

# 4D.5 (row 47) - Category + Description
Set-EmptyCellText $t.Rows.Item(47).Cells.Item(4) "Advanced"
Set-EmptyCellText $t.Rows.Item(47).Cells.Item(5) "Do Service Discovery"

# ---------------------------------------------------------------------------
# 5. Insert a new 4D.6 row after 4D.5, before the "05 (Debugging)" row.
# ---------------------------------------------------------------------------
$refRow = $t.Rows.Item(47)
$newRow = $t.Rows.Add($t.Rows.Item(48))
$newRow.Cells.Item(3).Range.Text = "4D.6"
$newRow.Cells.Item(4).Range.Text = "Advanced"
$newRow.Cells.Item(5).Range.Text = "Run the Advertising Scanner"
